# TC07_Canine_Filter_Breed-BrnMtnDog.xlsx - "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab Cypher query (cell B4 on the "startup" sheet) is corrected:
# the `File Type` and `Breed` columns are dropped from the RETURN clause
# (they duplicated/mismatched data already produced elsewhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesQuery = "MATCH (f:file)-->(parent)`n" +
    "WITH DISTINCT f, parent`n" +
    "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
    "WHERE demo.breed IN  ['Bernese Mountain Dog']`n" +
    "OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
    "OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "WITH DISTINCT f, parent, c, demo, diag, s`n" +
    "RETURN  coalesce(f.file_name, '') AS ``File Name``,`n" +
    "        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
    "        coalesce(f.file_description, '') AS ``Description``,`n" +
    "        coalesce(f.file_format, '') AS ``Format``,`n" +
    "        coalesce(f.file_size, '') AS ``Size``,`n" +
    "        coalesce(c.case_id, '') AS ``Case ID``,`n" +
    "        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
    "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Update the FilesTab query text in place.
$ws.Range("B4").Value = $newFilesQuery

# The row auto-sized shorter since the corrected query text is shorter.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection follows the edited cell.
[void]$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
